# Apply updated crypto price/volume data per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.737.03"
$ws.Range("E2").Value = "  -2.58%  "

# Row 3
$ws.Range("D3").Value = "3.268.37"
$ws.Range("E3").Value = "  -1.28%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.33"
$ws.Range("E5").Value = "  -1.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.76"
$ws.Range("E6").Value = "  -5.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  +3.58%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  -3.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  +0.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.397"
$ws.Range("E11").Value = "  -3.01%  "

# Row 12
$ws.Range("D12").Value = "3.838.44"
$ws.Range("E12").Value = "  -1.28%  "

# Row 13
$ws.Range("E13").Value = "  -3.85%  "

# Row 14
$ws.Range("D14").Value = "65.836.29"
$ws.Range("E14").Value = "  -2.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.34"
$ws.Range("E15").Value = "  -4.06%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.277.39"
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").Value = "  -3.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "435.37"
$ws.Range("E18").Value = "  -2.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.53"
$ws.Range("E19").Value = "  -2.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.10"
$ws.Range("E20").Value = "  -3.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.37"
$ws.Range("E21").Value = "  -5.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.24"
$ws.Range("E22").Value = "  -2.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("D24").Value = "3.424.54"
$ws.Range("E24").Value = "  -0.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.504"
$ws.Range("E25").Value = "  -2.53%  "

# Row 26
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.195"
$ws.Range("E26").Value = "  +3.81%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000112"
$ws.Range("E27").Value = "  -5.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.83"
$ws.Range("E28").Value = "  -2.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("E30").Value = "  -2.70%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.19"
$ws.Range("E31").Value = "  -3.29%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.11"
$ws.Range("E33").Value = "  -4.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.57"
$ws.Range("E34").Value = "  -3.53%  "

# Row 35
$ws.Range("E35").Value = "  -5.38%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.00"
$ws.Range("E36").Value = "  -1.78%  "

# Row 37
$ws.Range("E37").Value = "  -5.37%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.62"
$ws.Range("E38").Value = "  -2.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -4.45%  "

# Row 40
$ws.Range("D40").Value = "2.757.69"
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.771"
$ws.Range("E41").Value = "  -2.50%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.29"
$ws.Range("E42").Value = "  -4.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.22"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.01"
$ws.Range("E44").Value = "  -3.91%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0654"
$ws.Range("E45").Value = "  -2.75%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.27"
$ws.Range("E46").Value = "  -6.21%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "316.40"
$ws.Range("E47").Value = "  -3.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.24"
$ws.Range("E48").Value = "  -6.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0266"
$ws.Range("E49").Value = "  -2.81%  "

# Row 50
$ws.Range("E50").Value = "  +1.91%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.03%  "
